# Add a new "Italy" test-data sheet, modelled on the existing "Slovakia"
# sheet, and point the workbook at it as the active tab.

$wb = $excel.ActiveWorkbook
$slovakia = $wb.Worksheets.Item("Slovakia")

# Duplicate the Slovakia sheet (keeps layout, merges, col widths, styles,
# page setup, etc. identical) and drop the copy right after it.
$slovakia.Copy($null, $slovakia)
$italy = $wb.Worksheets.Item($slovakia.Index + 1)
$italy.Name = "Italy"

# Fill in the Italy-specific data. B4 (the product/user-story code) is set
# before B2 (the market name) so the new shared-string entries land in the
# same order as the target workbook.
$italy.Range("B4").Value = "NGC-3145/T2155"
$italy.Range("B2").Value = "Italy Market"

# B4 gets its own look: bigger, explicitly-black font with no cell border.
$italy.Range("B4").ClearFormats()
$italy.Range("B4").Font.Size = 12
$italy.Range("B4").Font.Color = 0

# Row 4 is made taller to fit the larger font.
$italy.Rows(4).RowHeight = 15.6

# Leave the new sheet selected at B2, and make it the active tab.
$italy.Activate()
$italy.Range("B2").Select()

# Slovakia goes back to an unremarkable "whole sheet" selection now that it
# is no longer the active tab.
$slovakia.Activate()
$slovakia.Cells.Select()
$italy.Activate()
